$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new test-case row: Red / 256 GB
$ws.Range("A5").Value = "Red"
$ws.Range("B5").Value = "256 GB"

# Move the active selection, matching the latest author view state
$ws.Range("C11").Select()
